{"js": "// Find the bullet paragraph \"Do not implement property generation (i.e. no cars\n// or street signs etc.)\" in the constraints list and insert a new bullet\n// \"Keep quality of models low.\" right after it, inheriting the same list\n// style/formatting (Word does this automatically for a paragraph inserted via\n// insertParagraph relative to an existing list paragraph).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Do not implement property generation\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the 'Do not implement property generation' paragraph.\");\n}\n\ntarget.insertParagraph(\"Keep quality of models low.\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Find the bullet paragraph \"Do not implement property generation (i.e. no cars\n# or street signs etc.)\" in the constraints list and insert a new bullet\n# \"Keep quality of models low.\" right after it, inheriting the same list\n# style/formatting from the paragraph it follows.\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Do not implement property generation*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the 'Do not implement property generation' paragraph.\"\n}\n\n$target.Range.InsertParagraphAfter()\n$newPara = $target.Next()\n$newPara.Range.Text = \"Keep quality of models low.\"\n"}
